$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell E1: "Mensajes" -> "Consulta de sensor"
$ws.Range("E1").Value = "Consulta de sensor"

# Update the active selection to E2 (as reflected in the sheetView)
$ws.Range("E2").Select()
